$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the previously used range clean (values + formats) so nothing
# from the old 2-row header / 6-row data layout leaks into the new
# 1-row header / 6-row data layout.
$ws.Range("A1:K8").Clear()

# ------------------------------------------------------------------
# Build a temporary named style so we get a cellXfs entry that has
# applyFont="1" but NOT applyNumberFormat="1" (fontId=1, numFmtId=0,
# xfId=0). Deleting the named style afterwards drops the
# cellStyleXfs/cellStyles bookkeeping but keeps the resulting xf in
# cellXfs - matching the target workbook's new style exactly.
# (NB: use IncludeNumber/Font here, NOT Styles(...).NumberFormat -
# the latter has been observed to silently drop the named style from
# the workbook's Styles collection.)
# ------------------------------------------------------------------
$headerStyle = $wb.Styles.Add("TmpHeaderStyle")
$headerStyle.IncludeNumber = $false
$headerStyle.Font.Name = "Arial"
$headerStyle.Font.Size = 9

# ------------------------------------------------------------------
# Row 1: new header row.
# ------------------------------------------------------------------
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"

$ws.Range("F1:K1").Style = "TmpHeaderStyle"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# Drop the helper named style now that its xf has been stamped onto
# F1:K1 - this removes the extra cellStyleXfs / cellStyles bookkeeping
# Excel creates for named styles, leaving only the one brand-new
# cellXfs entry the target workbook has.
$wb.Styles("TmpHeaderStyle").Delete()

# ------------------------------------------------------------------
# Rows 2..7: the six power-plant records (old rows 3..8 shifted up by
# one row), using the same per-column formats as the original sheet:
#   A,B,D,E -> Arial 9, "0"    number format
#   C       -> Arial 9, General
#   F..K    -> Arial 9, "0.00" number format
# Setting Font.Size/NumberFormat directly on the target ranges lets
# the engine dedupe onto the workbook's existing styles instead of
# minting new ones.
# ------------------------------------------------------------------
$data = @(
  @{ idx=1; idx2=108900; name="Wunderklingen"; start=1895; end=1968; f=5.5;   g=0.42;  h=0.41;               i=1.4;   j=1;     k=2.4 },
  @{ idx=2; idx2=106300; name="Engeweiher";    start=1909; end=1993; f=4;     g=5;     h=5;                  i=$null; j=$null; k=$null },
  @{ idx=3; idx2=108700; name="Eglisau";       start=1920; end=2012; f=500;   g=14.91; h=16.920000000000002; i=47.38; j=53.74; k=101.12 },
  @{ idx=4; idx2=106400; name="Neuhausen";     start=1951; end=2011; f=29.9;  g=2.8;   h=2.4500000000000002; i=10.35; j=10.4;  k=20.75 },
  @{ idx=5; idx2=106500; name="Rheinau";       start=1956; end=2005; f=400;   g=2.98;  h=2.92;               i=6.39;  j=13.19; k=19.579999999999998 },
  @{ idx=6; idx2=106200; name="Schaffhausen";  start=1964; end=$null; f=500;  g=22.57; h=19.84;              i=62.06; j=73.64; k=135.69999999999999 }
)

$row = 2
foreach ($rec in $data) {
  $ws.Range("A$row").Font.Size = 9
  $ws.Range("A$row").NumberFormat = "0"
  $ws.Range("A$row").Value = $rec.idx

  $ws.Range("B$row").Font.Size = 9
  $ws.Range("B$row").NumberFormat = "0"
  $ws.Range("B$row").Value = $rec.idx2

  $ws.Range("C$row").Font.Size = 9
  $ws.Range("C$row").Value = $rec.name

  $ws.Range("D$row").Font.Size = 9
  $ws.Range("D$row").NumberFormat = "0"
  $ws.Range("D$row").Value = $rec.start

  if ($null -ne $rec.end) {
    $ws.Range("E$row").Font.Size = 9
    $ws.Range("E$row").NumberFormat = "0"
    $ws.Range("E$row").Value = $rec.end
  }

  $ws.Range("F$row").Font.Size = 9
  $ws.Range("F$row").NumberFormat = "0.00"
  $ws.Range("F$row").Value = $rec.f

  $ws.Range("G$row").Font.Size = 9
  $ws.Range("G$row").NumberFormat = "0.00"
  $ws.Range("G$row").Value = $rec.g

  $ws.Range("H$row").Font.Size = 9
  $ws.Range("H$row").NumberFormat = "0.00"
  $ws.Range("H$row").Value = $rec.h

  if ($null -ne $rec.i) {
    $ws.Range("I$row").Font.Size = 9
    $ws.Range("I$row").NumberFormat = "0.00"
    $ws.Range("I$row").Value = $rec.i

    $ws.Range("J$row").Font.Size = 9
    $ws.Range("J$row").NumberFormat = "0.00"
    $ws.Range("J$row").Value = $rec.j

    $ws.Range("K$row").Font.Size = 9
    $ws.Range("K$row").NumberFormat = "0.00"
    $ws.Range("K$row").Value = $rec.k
  }

  $row = $row + 1
}

# ------------------------------------------------------------------
# Selection / dimension bookkeeping to match the edited workbook.
# ------------------------------------------------------------------
$ws.Range("A2:K2").Select()
